# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "51.701.59"
Set-TextValue $ws.Range("E2") "  +1.20%  "
Set-TextValue $ws.Range("D3") "3.063.65"
Set-TextValue $ws.Range("E3") "  +3.57%  "
Set-TextValue $ws.Range("D4") "0.998"
Set-TextValue $ws.Range("E4") "  -0.24%  "
Set-TextValue $ws.Range("D5") "385.15"
Set-TextValue $ws.Range("E5") "  +1.23%  "
Set-TextValue $ws.Range("D6") "103.35"
Set-TextValue $ws.Range("E6") "  +1.26%  "
Set-TextValue $ws.Range("D7") "0.545"
Set-TextValue $ws.Range("E7") "  +0.19%  "
Set-TextValue $ws.Range("E8") "  -0.01%  "
Set-TextValue $ws.Range("D9") "0.587"
Set-TextValue $ws.Range("E9") "  -0.54%  "
Set-TextValue $ws.Range("D10") "37.06"
Set-TextValue $ws.Range("E10") "  +1.79%  "
Set-TextValue $ws.Range("E11") "  +0.37%  "
Set-TextValue $ws.Range("E12") "  +1.44%  "
Set-TextValue $ws.Range("D13") "3.534.24"
Set-TextValue $ws.Range("E13") "  +3.22%  "
Set-TextValue $ws.Range("E14") "  +2.53%  "
Set-TextValue $ws.Range("D15") "7.76"
Set-TextValue $ws.Range("E15") "  -0.25%  "
Set-TextValue $ws.Range("D16") "3.054.04"
Set-TextValue $ws.Range("E16") "  +3.46%  "
Set-TextValue $ws.Range("D17") "0.984"
Set-TextValue $ws.Range("E17") "  -1.19%  "
Set-TextValue $ws.Range("D18") "10.58"
Set-TextValue $ws.Range("E18") "  -6.02%  "
Set-TextValue $ws.Range("D19") "51.744.99"
Set-TextValue $ws.Range("E19") "  +1.10%  "
Set-TextValue $ws.Range("D20") "3.12"
Set-TextValue $ws.Range("E20") "  -0.12%  "
Set-TextValue $ws.Range("D21") "12.50"
Set-TextValue $ws.Range("E21") "  +1.18%  "
Set-TextValue $ws.Range("D22") "0.0₃0967"
Set-TextValue $ws.Range("E22") "  +0.61%  "
Set-TextValue $ws.Range("D23") "70.13"
Set-TextValue $ws.Range("E23") "  -0.17%  "
Set-TextValue $ws.Range("D24") "267.88"
Set-TextValue $ws.Range("E24") "  +0.37%  "
Set-TextValue $ws.Range("E25") "  -3.96%  "
Set-TextValue $ws.Range("D26") "8.51"
Set-TextValue $ws.Range("E26") "  +8.09%  "
Set-TextValue $ws.Range("E27") "  +4.68%  "
Set-TextValue $ws.Range("D28") "7.36"
Set-TextValue $ws.Range("E28") "  +2.39%  "
Set-TextValue $ws.Range("D29") "26.73"
Set-TextValue $ws.Range("E29") "  +3.50%  "
Set-TextValue $ws.Range("E30") "  +0.14%  "
Set-TextValue $ws.Range("E31") "  -2.82%  "
Set-TextValue $ws.Range("D32") "10.32"
Set-TextValue $ws.Range("E32") "  +0.40%  "
Set-TextValue $ws.Range("D33") "34.29"
Set-TextValue $ws.Range("E33") "  -0.09%  "
Set-TextValue $ws.Range("E34") "  +0.79%  "
Set-TextValue $ws.Range("D35") "50.53"
Set-TextValue $ws.Range("E35") "  -1.04%  "
Set-TextValue $ws.Range("D36") "0.0445"
Set-TextValue $ws.Range("E36") "  +2.39%  "
Set-TextValue $ws.Range("E37") "  -0.21%  "
Set-TextValue $ws.Range("D38") "3.39"
Set-TextValue $ws.Range("E38") "  +5.28%  "
Set-TextValue $ws.Range("E39") "  +5.31%  "
Set-TextValue $ws.Range("D40") "17.13"
Set-TextValue $ws.Range("E40") "  +3.91%  "
Set-TextValue $ws.Range("D41") "1.88"
Set-TextValue $ws.Range("E41") "  +2.86%  "
Set-TextValue $ws.Range("D42") "128.32"
Set-TextValue $ws.Range("E42") "  +2.89%  "
Set-TextValue $ws.Range("E43") "  +0.02%  "
Set-TextValue $ws.Range("D44") "2.55"
Set-TextValue $ws.Range("E44") "  +1.71%  "
Set-TextValue $ws.Range("D45") "3.70"
Set-TextValue $ws.Range("E45") "  +5.11%  "
Set-TextValue $ws.Range("D46") "22.04"
Set-TextValue $ws.Range("E46") "  +3.27%  "
Set-TextValue $ws.Range("D47") "2.53"
Set-TextValue $ws.Range("E47") "  +6.70%  "
Set-TextValue $ws.Range("E48") "  +3.22%  "
Set-TextValue $ws.Range("D49") "2.043.36"
Set-TextValue $ws.Range("E49") "  -0.09%  "
Set-TextValue $ws.Range("D50") "3.358.10"
Set-TextValue $ws.Range("E50") "  +3.27%  "
Set-TextValue $ws.Range("E51") "  +7.70%  "
